# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback sync
# completed successfully for both target locales (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (shown on the Overview sheet per-locale, and on each locale sheet's
#     Status column).
#   - Each locale sheet's "Latest Handback DateTime" is refreshed to the
#     timestamp of this run.
#   - The previous "handback file is not the latest" Error Detail is
#     cleared now that everything is back in sync.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus                 # Status
$zhcn.Range("K2").Value = "2016-08-23 22:46:31"       # Latest Handback DateTime
# Error Detail is cleared but stays a (blank) text cell rather than being
# removed outright, so a leading quote forces text-empty instead of a
# numeric/blank cell, then the style is reset back to Normal.
$zhcn.Range("P2").Value = "'"
$zhcn.Range("P2").Style = "Normal"
$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value = $newStatus                 # Status
$dede.Range("K2").Value = "2016-08-23 22:46:38"       # Latest Handback DateTime
$dede.Range("P2").Value = "'"
$dede.Range("P2").Style = "Normal"
$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
